# Insert a new row at position 56, shifting existing rows 56-150 down to 57-151.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("56:56").Insert()

# Populate the newly inserted row 56 with its data.
$ws.Cells.Item(56, 1).Value = 4
$ws.Cells.Item(56, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(56, 3).Value = "Los Lagos"
$ws.Cells.Item(56, 4).Value = 44544
$ws.Cells.Item(56, 5).Value = 10
$ws.Cells.Item(56, 6).Value = 100112028
$ws.Cells.Item(56, 7).Value = "Sandia"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Tercera"
$ws.Cells.Item(56, 10).Value = 800
$ws.Cells.Item(56, 11).Value = 2500
$ws.Cells.Item(56, 12).Value = 2500
$ws.Cells.Item(56, 13).Value = 2500
$ws.Cells.Item(56, 14).Value = "$/unidad"
$ws.Cells.Item(56, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(56, 16).Value = 2500
$ws.Cells.Item(56, 17).Value = 1
$ws.Cells.Item(56, 18).Value = "Hortaliza"
